$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize every data row (1-88) to the new row height used throughout
# the sheet, then mark a specific subset of rows as hidden (filtered
# view), matching the target workbook state.
$ws.Range("A1:C88").RowHeight = 17.25

$hiddenRows = @(1, 10, 12, 29, 32, 35, 40, 41, 59, 64, 71, 75, 80, 88)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}
